# Update the "Förändrad" (Changed) date column (C) for rows 2-46
# from serial date 45188 (2023-09-19) to 45189 (2023-09-20).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 46; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    $cell.Value = 45189
}
